# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF.
# Copy formatting from an existing header cell (e.g. AC1) so the new headers
# match the style (bold, bordered, centered) used by the rest of row 1, then
# set the text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-52: season record values (Wins=76, Losses=86, Ties=0)
$lastRow = 52
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 76
    $ws.Cells.Item($r, 31).Value = 86
    $ws.Cells.Item($r, 32).Value = 0
}
